# Edit script: add Source 5/6/7 (NYC, UK Public Health Profiles, Houston) data
# to "profile assessment (transposed)" sheet, drop the autofilter on Table1,
# and tidy up selection / column widths to match the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Table1 no longer shows the header-row autofilter dropdowns.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.ShowAutoFilter = $false

# ---------------------------------------------------------------------------
# 2. New source footnote rows at the bottom of the sheet.
#    Entered in the same order the author appears to have typed them:
#    the Source8/9/10 placeholders first, then each new source label
#    immediately before its data column.
# ---------------------------------------------------------------------------
$ws.Range("A29").Value = "Source8= "
$ws.Range("A30").Value = "Source9= "
$ws.Range("A31").Value = "Source10+"

$ws.Range("A26").Value = "Source5= NYC.gov Community Health Profiles"

# ---------------------------------------------------------------------------
# 3. Column G = Source 5 (NYC.gov Community Health Profiles)
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = "59 communites"
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 22
$ws.Range("G5").Value = 7
$ws.Range("G6").Value = 8
$ws.Range("G7").Value = "Yes "
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = "Yes "
$ws.Range("G11").Value = 0
$ws.Range("G12").Value = 20
$ws.Range("G14").Value = "Yes "
$ws.Range("G15").Value = "Yes "
$ws.Range("G16").Value = "Yes "
$ws.Range("G17").Value = "Yes "
$ws.Range("G18").Value = "Yes-languages & font sizing"

# ---------------------------------------------------------------------------
# 4. Column H = Source 6 (Department of Health & Social Care Public Health
#    Profiles, United Kingdom)
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "Source6= Department of Health & Social Care Public Health Profiles (United Kingdom)"

$ws.Range("H2").Value = "4 Area Types"
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = "Yes"
$ws.Range("H5").Value = "Yes"
$ws.Range("H6").Value = "7/38"
$ws.Range("H7").Value = "Yes"
$ws.Range("H8").Value = "Yes"
$ws.Range("H9").Value = "Yes"
$ws.Range("H10").Value = "Yes"
$ws.Range("H11").Value = 0
$ws.Range("H13").Value = "interactive web page and PDF available"
$ws.Range("H14").Value = "Yes"
$ws.Range("H15").Value = "No"
$ws.Range("H16").Value = "Yes"
$ws.Range("H17").Value = "Yes"
$ws.Range("H18").Value = "No"

# ---------------------------------------------------------------------------
# 5. Column I = Source 7 (Houston Health Department Community Health
#    Profiles and Reports)
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Source7= Houston Health Department Community Health Profiles and Reports"

$ws.Range("I2").Value = "4 Communities"
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 22
$ws.Range("I5").Value = 0
$ws.Range("I6").Value = "3/20"
$ws.Range("I7").Value = "Yes"
$ws.Range("I8").Value = 4
$ws.Range("I9").Value = "No"
$ws.Range("I10").Value = "Yes"
$ws.Range("I11").Value = 0
$ws.Range("I12").Value = 31
$ws.Range("I14").Value = "Yes"
$ws.Range("I15").Value = "Yes"
$ws.Range("I16").Value = "Yes-Houston State of Health, CDC, Healthy People 2020"
$ws.Range("I17").Value = "No"
$ws.Range("I18").Value = "No"

# ---------------------------------------------------------------------------
# 6. H6 / I6 hold "7/38" and "3/20" -- Excel stores these as right-aligned
#    text (numFmt "@") rather than dates/fractions. Format H6 directly, then
#    copy/paste-format onto I6 so both cells share a single new style entry.
# ---------------------------------------------------------------------------
$h6 = $ws.Range("H6")
$h6.NumberFormat = "@"
$h6.HorizontalAlignment = -4152
$h6.Copy()
$ws.Range("I6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 7. Column widths for C:L were re-fit after the new columns were populated.
# ---------------------------------------------------------------------------
$ws.Range("C1").EntireColumn.ColumnWidth = 10.5
$ws.Range("D1").EntireColumn.ColumnWidth = 10.33
$ws.Range("E1").EntireColumn.ColumnWidth = 9.83
$ws.Range("F1").EntireColumn.ColumnWidth = 10
$ws.Range("G1").EntireColumn.ColumnWidth = 10.67
$ws.Range("H1").EntireColumn.ColumnWidth = 10.33
$ws.Range("I1").EntireColumn.ColumnWidth = 10.17
$ws.Range("J1").EntireColumn.ColumnWidth = 10.83
$ws.Range("K1").EntireColumn.ColumnWidth = 10.67
$ws.Range("L1").EntireColumn.ColumnWidth = 11.33

# ---------------------------------------------------------------------------
# 8. Final selection / view state.
# ---------------------------------------------------------------------------
$ws.Range("J2").Select()
